$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.918.06"
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = "'2.530.49"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'305.22"
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").Value = "'101.90"
$ws.Range("E6").Value = '  +7.96%  '
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").Value = "'37.82"
$ws.Range("E10").Value = '  +4.04%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = "'7.61"
$ws.Range("E12").Value = '  -1.75%  '
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = "'2.924.19"
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = "'2.503.01"
$ws.Range("E15").Value = '  -1.77%  '
$ws.Range("D16").Value = "'15.18"
$ws.Range("E16").Value = '  +6.96%  '
$ws.Range("D17").Value = "'0.865"
$ws.Range("E17").Value = '  -1.53%  '
$ws.Range("D18").Value = "'42.925.04"
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = "'13.16"
$ws.Range("E19").Value = '  +3.67%  '
$ws.Range("D20").Value = "'0.0₃0985"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = "'6.50"
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = "'252.65"
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  -2.83%  '
$ws.Range("D26").Value = "'27.19"
$ws.Range("E26").Value = '  -6.04%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("E28").Value = '  +8.76%  '
$ws.Range("D29").Value = "'10.34"
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").Value = "'39.10"
$ws.Range("E30").Value = '  +5.37%  '
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").Value = "'157.75"
$ws.Range("E32").Value = '  +2.33%  '
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("E36").Value = '  -3.32%  '
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("D39").Value = "'24.24"
$ws.Range("E39").Value = '  +5.65%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +2.67%  '
$ws.Range("D42").Value = "'2.10"
$ws.Range("E42").Value = '  -3.23%  '
$ws.Range("D43").Value = "'3.89"
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = "'0.0304"
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = "'2.042.68"
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("D47").Value = "'86.47"
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("E48").Value = '  -3.69%  '
$ws.Range("D49").Value = "'2.781.62"
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("E50").Value = '  +1.86%  '
$ws.Range("D51").Value = "'103.03"
$ws.Range("E51").Value = '  -2.76%  '
